# Refresh the "cryptos" price/volume snapshot (Price = column D, Volume(1h) = column E).
# Values are written as literal text (matching the source data, which stores prices/
# percentages as text, e.g. "2.644.42" or "  -3.45%  ") rather than as numbers.
# A leading apostrophe is used for D-column values that would otherwise be
# auto-interpreted by Excel as a number (so trailing zeros / exact text are preserved).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "67.556.16"
$ws.Range("E2").Value = "  -2.21%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.644.42"
$ws.Range("E3").Value = "  -3.45%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5: BNB
$ws.Range("D5").Value = "'598.55"
$ws.Range("E5").Value = "  -1.35%  "

# Row 6: Solana
$ws.Range("D6").Value = "'168.05"
$ws.Range("E6").Value = "  -1.24%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8: XRP
$ws.Range("E8").Value = "  -0.67%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "2.643.87"
$ws.Range("E9").Value = "  -3.43%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  -1.60%  "

# Row 11: TRON
$ws.Range("E11").Value = "  +1.95%  "

# Row 12: Cardano
$ws.Range("E12").Value = "  -1.26%  "

# Row 13: Toncoin
$ws.Range("E13").Value = "  -2.21%  "

# Row 14: Avalanche
$ws.Range("D14").Value = "'28.09"
$ws.Range("E14").Value = "  -2.72%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.118.85"
$ws.Range("E15").Value = "  -3.65%  "

# Row 16: ShibaInu
$ws.Range("E16").Value = "  -3.81%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "67.677.45"
$ws.Range("E17").Value = "  -1.92%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "2.641.10"
$ws.Range("E18").Value = "  -2.56%  "

# Row 19: Chainlink
$ws.Range("D19").Value = "'11.92"
$ws.Range("E19").Value = "  -0.20%  "

# Row 20: Uniswap
$ws.Range("D20").Value = "'7.92"
$ws.Range("E20").Value = "  +2.63%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'363.45"
$ws.Range("E21").Value = "  -3.50%  "

# Row 22: Polkadot
$ws.Range("E22").Value = "  -3.28%  "

# Row 23: NEARProtocol
$ws.Range("D23").Value = "'4.80"
$ws.Range("E23").Value = "  -4.21%  "

# Row 24: Aptos
$ws.Range("D24").Value = "'11.05"
$ws.Range("E24").Value = "  +8.67%  "

# Row 25: SuiNetwork
$ws.Range("D25").Value = "'2.02"
$ws.Range("E25").Value = "  -5.24%  "

# Row 26: Dai
$ws.Range("E26").Value = "  -0.03%  "

# Row 27: Litecoin
$ws.Range("D27").Value = "'70.76"
$ws.Range("E27").Value = "  -4.39%  "

# Row 28: WrappedeETH
$ws.Range("D28").Value = "2.776.15"
$ws.Range("E28").Value = "  -3.35%  "

# Row 29: PEPE
$ws.Range("E29").Value = "  -3.87%  "

# Row 30: Binance-PegBSC-USD
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.16%  "

# Row 31: Bittensor
$ws.Range("D31").Value = "'557.34"
$ws.Range("E31").Value = "  -5.66%  "

# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").Value = "'8.04"
$ws.Range("E32").Value = "  -4.11%  "

# Row 33: Fetch.AI
$ws.Range("E33").Value = "  -4.25%  "

# Row 34: PancakeSwap
$ws.Range("E34").Value = "  -2.55%  "

# Row 35: Kaspa
$ws.Range("E35").Value = "  -0.03%  "

# Row 36: FirstDigitalUSD
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.02%  "

# Row 37: ImmutableX
$ws.Range("E37").Value = "  -5.03%  "

# Row 38: Monero
$ws.Range("D38").Value = "'157.97"
$ws.Range("E38").Value = "  -2.72%  "

# Row 39: EthereumClassic
$ws.Range("D39").Value = "'19.42"
$ws.Range("E39").Value = "  -3.34%  "

# Row 40: PolygonEcosystemToken
$ws.Range("E40").Value = "  -2.71%  "

# Row 41: RenderToken
$ws.Range("D41").Value = "'5.30"
$ws.Range("E41").Value = "  -4.02%  "

# Row 42: Stacks
$ws.Range("D42").Value = "'1.84"
$ws.Range("E42").Value = "  -4.87%  "

# Row 43: WhiteBITCoin
$ws.Range("D43").Value = "'17.94"
$ws.Range("E43").Value = "  -0.45%  "

# Row 44: dogwifhat
$ws.Range("E44").Value = "  -6.20%  "

# Row 45: USDe
$ws.Range("E45").Value = "  +0.00%  "

# Row 46: OKB
$ws.Range("E46").Value = "  -2.41%  "

# Row 47: BabyDogeCoin
$ws.Range("E47").Value = "  -3.57%  "

# Row 48: ARBITRUM
$ws.Range("D48").Value = "'0.597"
$ws.Range("E48").Value = "  -1.78%  "

# Row 49: Aave
$ws.Range("D49").Value = "'153.94"
$ws.Range("E49").Value = "  -1.60%  "

# Row 50: Filecoin
$ws.Range("E50").Value = "  -2.35%  "

# Row 51: Optimism
$ws.Range("E51").Value = "  -4.15%  "
